$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Column order for the per-row data block (D through AJ).
$cols = @("D","E","F","G","H","I","J","K","L","M","N","O","P","Q","R","S","T","U","V","W","X","Y","Z","AA","AB","AC","AD","AE","AF","AG","AH","AI","AJ")

# Rows 2-6 (2014/12 .. 2018/12 IFRS-consolidated actuals): every metric in
# D:AJ is refreshed with the corrected figures, and the "FCF" column (U) is
# dropped (blanked out) for all of them.
$rows = @{
    2 = @(3688,207,207,207,154,154,0,8958,4985,3973,3971,3,668,-1978,-191,2167,5,$null,717,5.6,4.17,3.98,1.98,125.46,543.64,1150,13.17,40475,0.37,1000,6.6,64.77,10369886)
    3 = @(7089,302,302,328,248,248,0,12078,7837,4240,4226,15,668,-616,98,834,8,$null,728,4.26,3.5,6.05,2.36,184.82,583.57,1855,8.63,43071,0.37,1200,7.5,48.07,10369886)
    4 = @(7537,353,353,374,278,278,0,14334,9751,4583,4570,13,668,1383,-169,-724,16,$null,691,4.68,3.68,6.31,2.1,212.75,634.88,2077,9.05,46583,0.4,1200,6.38,42.93,10369886)
    5 = @(5663,473,473,499,376,374,2,13888,9095,4793,4780,13,668,-2592,-34,2514,6,$null,1901,8.35,6.64,8,2.65,189.75,666.27,2797,8.99,48723,0.52,1200,4.77,31.87,10369886)
    6 = @(6272,373,373,397,282,290,$null,14351,9702,4649,4617,$null,668,5,23,-67,8,$null,3048,5.95,4.5,6.17,2.05,208.69,682.91,2170,10.34,51862,0.43,1200,5.35,37.33,10369886)
}

foreach ($r in $rows.Keys) {
    $vals = $rows[$r]
    for ($i = 0; $i -lt $cols.Length; $i++) {
        $ws.Range("$($cols[$i])$r").Value = $vals[$i]
    }
}

# Rows 7-9 (2019/12(E) .. 2021/12(E) forecasts): these columns are removed
# entirely, leaving only the row number / period / label columns (A:C).
$ws.Range("D7:AJ9").ClearContents()
